# vacinas.xlsx - "finalizando tarefa de vacinas FPOO"
#
# Adds a small "scratch" side-table in J12:M13 (copies of the first two
# vaccine records used while testing), and finalises the "alterar o 2
# elemento da lista" edit by overwriting row 15 (list index 1) with the
# Pamonha/Carlos/doencas/2022-08-10 record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12: Rodolfo / Rosana / gripe / 10-28-2022 ---------------------
$ws.Range("J12").Value = "Rodolfo"
$ws.Range("K12").Value = "Rosana"
$ws.Range("L12").Value = "gripe"
$ws.Range("M12").Value = 44862
$ws.Range("M12").NumberFormat = "mm-dd-yy"

# --- Extra data alongside row 13 (existing "Excluir o elemento 1" header) --
# Rodrigo / James / sarna / 10-30-2022
$ws.Range("J13").Value = "Rodrigo"
$ws.Range("K13").Value = "James"
$ws.Range("L13").Value = "sarna"
$ws.Range("M13").Value = 44864

# Reuse M12's date number formatting for M13 (keeps a single shared style)
$ws.Range("M12").Copy()
$ws.Range("M13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 15 (list index 1): Pudim/James/sarna -> Pamonha/Carlos/doencas ----
$ws.Range("B15").Value = "Pamonha"
$ws.Range("C15").Value = "Carlos"
$ws.Range("D15").Value = "doencas"
$ws.Range("E15").Value = 44783

# --- Column sizing for the new L/M columns ----------------------------------
$ws.Columns.Item(12).ColumnWidth = 4.75
$ws.Columns.Item(13).ColumnWidth = 10.25

# --- View zoom ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 145
